$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "ffff341494fd-1949-4346-98a4-80679563b85d.md"
$ws.Range("B2").Value = "e2e\ffff341494fd-1949-4346-98a4-80679563b85d.md"
$ws.Range("G2").Value = "2016-09-05 23:16:38"
$ws.Range("A3").Value = "ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md"
$ws.Range("B3").Value = "e2e\ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md"
$ws.Range("A4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("B4").Value = "e2e\45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-09-05 23:20:23"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') { $h.TextToDisplay = "e2e\ffff341494fd-1949-4346-98a4-80679563b85d.md" }
    if ($addr -eq '$B$3') { $h.TextToDisplay = "e2e\ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md" }
    if ($addr -eq '$B$4') { $h.TextToDisplay = "e2e\45fba168-d119-4e17-adff-de1f86fcaaa0.md" }
}

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "ffff341494fd-1949-4346-98a4-80679563b85d.md"
$ws.Range("G2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.5699311b2f211a1627904973773dd11918028c8f.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-05 23:16:32"
$ws.Range("I2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.md"
$ws.Range("J2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.5699311b2f211a1627904973773dd11918028c8f.zh-cn.xlf"
$ws.Range("K2").Value = "2016-09-05 23:16:50"
$ws.Range("A3").Value = "ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md"
$ws.Range("F3").Value = "'True"
$ws.Range("A4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.abab97d577e263667c0f194b6704ead82a7ff137.zh-cn.xlf"
$ws.Range("H4").Value = "2016-09-05 23:20:18"
$ws.Range("I4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("J4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.abab97d577e263667c0f194b6704ead82a7ff137.zh-cn.xlf"
$ws.Range("K4").Value = "2016-09-05 23:19:43"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/618969977d7595762a5e3db2e7cc993ff102639f/e2e/45fba168-d119-4e17-adff-de1f86fcaaa0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d32394f0ccb1f006630bc413435d47a5ec4e17a1/e2e/45fba168-d119-4e17-adff-de1f86fcaaa0.md."

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff341494fd-1949-4346-98a4-80679563b85d.md" }
    if ($addr -eq '$A$3') { $h.TextToDisplay = "ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md" }
    if ($addr -eq '$A$4') { $h.TextToDisplay = "45fba168-d119-4e17-adff-de1f86fcaaa0.md" }
    if ($addr -eq '$I$2') { $h.TextToDisplay = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.md" }
    if ($addr -eq '$I$4') { $h.TextToDisplay = "45fba168-d119-4e17-adff-de1f86fcaaa0.md" }
}

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "ffff341494fd-1949-4346-98a4-80679563b85d.md"
$ws.Range("G2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.5699311b2f211a1627904973773dd11918028c8f.de-de.xlf"
$ws.Range("H2").Value = "2016-09-05 23:16:38"
$ws.Range("I2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.md"
$ws.Range("J2").Value = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.5699311b2f211a1627904973773dd11918028c8f.de-de.xlf"
$ws.Range("K2").Value = "2016-09-05 23:16:58"
$ws.Range("A3").Value = "ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md"
$ws.Range("F3").Value = "'True"
$ws.Range("A4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.abab97d577e263667c0f194b6704ead82a7ff137.de-de.xlf"
$ws.Range("H4").Value = "2016-09-05 23:20:23"
$ws.Range("I4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.md"
$ws.Range("J4").Value = "45fba168-d119-4e17-adff-de1f86fcaaa0.abab97d577e263667c0f194b6704ead82a7ff137.de-de.xlf"
$ws.Range("K4").Value = "2016-09-05 23:19:51"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/618969977d7595762a5e3db2e7cc993ff102639f/e2e/45fba168-d119-4e17-adff-de1f86fcaaa0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d32394f0ccb1f006630bc413435d47a5ec4e17a1/e2e/45fba168-d119-4e17-adff-de1f86fcaaa0.md."

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff341494fd-1949-4346-98a4-80679563b85d.md" }
    if ($addr -eq '$A$3') { $h.TextToDisplay = "ffffff6407df4d-3ba0-4f1e-8a41-80dbc6736554.md" }
    if ($addr -eq '$A$4') { $h.TextToDisplay = "45fba168-d119-4e17-adff-de1f86fcaaa0.md" }
    if ($addr -eq '$I$2') { $h.TextToDisplay = "73a0d2a1-89a1-481a-ab4c-4ce0e12a45c0.md" }
    if ($addr -eq '$I$4') { $h.TextToDisplay = "45fba168-d119-4e17-adff-de1f86fcaaa0.md" }
}

# ---- Column width updates ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(16).ColumnWidth = 39.16666666666666
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(16).ColumnWidth = 39.16666666666666
